$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "29.083.98"
$ws.Cells.Item(2, 5).Value = "  +0.61%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.832.12"
$ws.Cells.Item(3, 5).Value = "  +0.01%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  +0.31%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "238.17"
$ws.Cells.Item(5, 5).Value = "  -2.74%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.6823"
$ws.Cells.Item(6, 5).Value = "  -1.35%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.003"
$ws.Cells.Item(7, 5).Value = "  +0.33%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3008"
$ws.Cells.Item(8, 5).Value = "  -1.00%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07458"
$ws.Cells.Item(9, 5).Value = "  -2.70%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "23.07"
$ws.Cells.Item(10, 5).Value = "  -0.94%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07676"

# Row 12
$ws.Cells.Item(12, 4).Value = "1.832.54"
$ws.Cells.Item(12, 5).Value = "  +0.46%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.055"
$ws.Cells.Item(13, 5).Value = "  -0.65%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.6816"
$ws.Cells.Item(14, 5).Value = "  +0.10%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "87.24"
$ws.Cells.Item(15, 5).Value = "  -6.20%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "6.236"
$ws.Cells.Item(16, 5).Value = "  -4.55%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "29.100.04"
$ws.Cells.Item(17, 5).Value = "  +0.60%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000008162"
$ws.Cells.Item(18, 5).Value = "  -0.87%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "2.089.30"
$ws.Cells.Item(19, 5).Value = "  +0.69%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "227.23"
$ws.Cells.Item(20, 5).Value = "  -5.36%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "12.52"
$ws.Cells.Item(21, 5).Value = "  -1.11%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.002"
$ws.Cells.Item(22, 5).Value = "  +0.22%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.412"
$ws.Cells.Item(23, 5).Value = "  -0.59%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.003"
$ws.Cells.Item(24, 5).Value = "  +0.33%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.1453"
$ws.Cells.Item(25, 5).Value = "  -2.91%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "159.88"
$ws.Cells.Item(26, 5).Value = "  +1.03%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.764"
$ws.Cells.Item(27, 5).Value = "  +0.38%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "18.06"
$ws.Cells.Item(28, 5).Value = "  -0.42%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.496"
$ws.Cells.Item(29, 5).Value = "  -2.76%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.250"
$ws.Cells.Item(30, 5).Value = "  +0.74%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.140"
$ws.Cells.Item(31, 5).Value = "  +0.11%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.204"
$ws.Cells.Item(32, 5).Value = "  +0.78%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05138"
$ws.Cells.Item(33, 5).Value = "  +0.59%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.7708"
$ws.Cells.Item(34, 5).Value = "  -0.71%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.829"
$ws.Cells.Item(35, 5).Value = "  -1.22%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.132"
$ws.Cells.Item(36, 5).Value = "  -0.88%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.680"
$ws.Cells.Item(37, 5).Value = "  -0.40%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "1.306.86"
$ws.Cells.Item(38, 5).Value = "  +2.16%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01836"
$ws.Cells.Item(39, 5).Value = "  -1.05%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.691"
$ws.Cells.Item(40, 5).Value = "  -0.25%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.9351"
$ws.Cells.Item(41, 5).Value = "  -2.04%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.783"
$ws.Cells.Item(42, 5).Value = "  -5.72%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Quant"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "104.60"
$ws.Cells.Item(43, 5).Value = "  -2.03%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.002"
$ws.Cells.Item(44, 5).Value = "  +0.27%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "RocketPoolETH"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(45, 4).Value = "1.987.87"
$ws.Cells.Item(45, 5).Value = "  +0.65%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "65.25"
$ws.Cells.Item(46, 5).Value = "  +2.48%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Mantle"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5196"
$ws.Cells.Item(47, 5).Value = "  +0.58%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.22%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.543"
$ws.Cells.Item(49, 5).Value = "  -1.35%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.762"
$ws.Cells.Item(50, 5).Value = "  +0.65%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "XinFinNetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.07343"
$ws.Cells.Item(51, 5).Value = "  +20.96%  "
